# Update crypto price/volume figures (inline text cells) per the latest GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '41.813.48'
$cell.Style = 'Normal'
$cell = $ws.Range('E2')
$cell.NumberFormat = '@'
$cell.Value = '  +0.69%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '2.478.78'
$cell.Style = 'Normal'
$cell = $ws.Range('E3')
$cell.NumberFormat = '@'
$cell.Value = '  +0.37%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E4')
$cell.NumberFormat = '@'
$cell.Value = '  +0.05%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '319.17'
$cell.Style = 'Normal'
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '93.30'
$cell.Style = 'Normal'
$cell = $ws.Range('E6')
$cell.NumberFormat = '@'
$cell.Value = '  +1.46%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E7')
$cell.NumberFormat = '@'
$cell.Value = '  +0.43%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E8')
$cell.NumberFormat = '@'
$cell.Value = '  +0.11%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E9')
$cell.NumberFormat = '@'
$cell.Value = '  +0.30%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '0.0878'
$cell.Style = 'Normal'
$cell = $ws.Range('E10')
$cell.NumberFormat = '@'
$cell.Value = '  +10.81%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '33.36'
$cell.Style = 'Normal'
$cell = $ws.Range('E11')
$cell.NumberFormat = '@'
$cell.Value = '  +2.59%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E12')
$cell.NumberFormat = '@'
$cell.Value = '  +0.74%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '2.860.56'
$cell.Style = 'Normal'
$cell = $ws.Range('E13')
$cell.NumberFormat = '@'
$cell.Value = '  +0.37%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '6.93'
$cell.Style = 'Normal'
$cell = $ws.Range('E14')
$cell.NumberFormat = '@'
$cell.Value = '  +0.98%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '15.75'
$cell.Style = 'Normal'
$cell = $ws.Range('E15')
$cell.NumberFormat = '@'
$cell.Value = '  -1.34%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '2.483.05'
$cell.Style = 'Normal'
$cell = $ws.Range('E16')
$cell.NumberFormat = '@'
$cell.Value = '  +0.47%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '0.799'
$cell.Style = 'Normal'
$cell = $ws.Range('E17')
$cell.NumberFormat = '@'
$cell.Value = '  +2.71%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '41.764.11'
$cell.Style = 'Normal'
$cell = $ws.Range('E18')
$cell.NumberFormat = '@'
$cell.Value = '  +0.53%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '6.47'
$cell.Style = 'Normal'
$cell = $ws.Range('E19')
$cell.NumberFormat = '@'
$cell.Value = '  -0.66%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E20')
$cell.NumberFormat = '@'
$cell.Value = '  +0.71%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '71.22'
$cell.Style = 'Normal'
$cell = $ws.Range('E21')
$cell.NumberFormat = '@'
$cell.Value = '  +0.29%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '11.35'
$cell.Style = 'Normal'
$cell = $ws.Range('E22')
$cell.NumberFormat = '@'
$cell.Value = '  +1.99%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '240.57'
$cell.Style = 'Normal'
$cell = $ws.Range('E23')
$cell.NumberFormat = '@'
$cell.Value = '  +1.09%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E24')
$cell.NumberFormat = '@'
$cell.Value = '  +1.23%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E25')
$cell.NumberFormat = '@'
$cell.Value = '  +2.25%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E26')
$cell.NumberFormat = '@'
$cell.Value = '  +0.00%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '24.81'
$cell.Style = 'Normal'
$cell = $ws.Range('E27')
$cell.NumberFormat = '@'
$cell.Value = '  +0.18%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E28')
$cell.NumberFormat = '@'
$cell.Value = '  +1.04%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '9.80'
$cell.Style = 'Normal'
$cell = $ws.Range('E29')
$cell.NumberFormat = '@'
$cell.Value = '  +0.89%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '36.87'
$cell.Style = 'Normal'
$cell = $ws.Range('E30')
$cell.NumberFormat = '@'
$cell.Value = '  +3.89%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '157.89'
$cell.Style = 'Normal'
$cell = $ws.Range('E31')
$cell.NumberFormat = '@'
$cell.Value = '  +1.22%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '5.52'
$cell.Style = 'Normal'
$cell = $ws.Range('E32')
$cell.NumberFormat = '@'
$cell.Value = '  +1.21%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E33')
$cell.NumberFormat = '@'
$cell.Value = '  -0.04%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '0.0767'
$cell.Style = 'Normal'
$cell = $ws.Range('E34')
$cell.NumberFormat = '@'
$cell.Value = '  +0.84%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E35')
$cell.NumberFormat = '@'
$cell.Value = '  +0.09%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E36')
$cell.NumberFormat = '@'
$cell.Value = '  +1.53%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '1.88'
$cell.Style = 'Normal'
$cell = $ws.Range('E37')
$cell.NumberFormat = '@'
$cell.Value = '  +4.70%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '2.94'
$cell.Style = 'Normal'
$cell = $ws.Range('E38')
$cell.NumberFormat = '@'
$cell.Value = '  +1.50%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E39')
$cell.NumberFormat = '@'
$cell.Value = '  +1.69%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E40')
$cell.NumberFormat = '@'
$cell.Value = '  +0.73%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '2.55'
$cell.Style = 'Normal'
$cell = $ws.Range('E41')
$cell.NumberFormat = '@'
$cell.Value = '  +9.06%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E42')
$cell.NumberFormat = '@'
$cell.Value = '  +0.34%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '2.002.29'
$cell.Style = 'Normal'
$cell = $ws.Range('E43')
$cell.NumberFormat = '@'
$cell.Value = '  +2.86%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '19.22'
$cell.Style = 'Normal'
$cell = $ws.Range('E44')
$cell.NumberFormat = '@'
$cell.Value = '  +1.73%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '0.0286'
$cell.Style = 'Normal'
$cell = $ws.Range('E45')
$cell.NumberFormat = '@'
$cell.Value = '  +0.75%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E46')
$cell.NumberFormat = '@'
$cell.Value = '  +3.13%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '9.49'
$cell.Style = 'Normal'
$cell = $ws.Range('E47')
$cell.NumberFormat = '@'
$cell.Value = '  +4.16%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '2.715.10'
$cell.Style = 'Normal'
$cell = $ws.Range('E48')
$cell.NumberFormat = '@'
$cell.Value = '  +0.22%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '98.56'
$cell.Style = 'Normal'
$cell = $ws.Range('E49')
$cell.NumberFormat = '@'
$cell.Value = '  +1.19%  '
$cell.Style = 'Normal'
$cell = $ws.Range('E50')
$cell.NumberFormat = '@'
$cell.Value = '  +4.55%  '
$cell.Style = 'Normal'
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '67.43'
$cell.Style = 'Normal'
$cell = $ws.Range('E51')
$cell.NumberFormat = '@'
$cell.Value = '  +0.24%  '
$cell.Style = 'Normal'
